$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (pushes old rows 13-16 down to 14-17)
$ws.Rows.Item(13).Insert()

# Rename "Antenna" (row 12) label to "Antenna Type"
$ws.Range("A12").Value = "Antenna Type"

# Fill new row 13: "Antenna Beam Aperture (degrees)" with per-site values
$ws.Range("A13").Value = "Antenna Beam Aperture (degrees)"
$ws.Range("A13").Style = $ws.Range("A12").Style

$values13 = @(70, 80, 55, 80, 70, 70, 70, 70, 70, 70, 70, 70)
for ($i = 0; $i -lt $values13.Length; $i++) {
    $col = 2 + $i
    $ws.Cells.Item(13, $col).Value = $values13[$i]
}

# Fix the transmitter frequency row (now row 15) values: 14000000000 -> 1400000000
for ($col = 2; $col -le 13; $col++) {
    $ws.Cells.Item(15, $col).Value = 1400000000
}

# Update the selection shown in the sheet view
$ws.Range("A2:A16").Select()

Write-Host "Done"
